$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1) for new columns I and J, mirroring style/format of existing headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting from the neighboring header cell (H1) onto I1:J1 so the
# new header cells share the same bold / bordered / centered style.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data values for columns I and J, rows 2-16
$data = @{
    2  = @(7, 8)
    3  = @(9, 9)
    4  = @(8, 8)
    5  = @(8, 8)
    6  = @(5, 8)
    7  = @(6, 8)
    8  = @(8, 8)
    9  = @(9, 9)
    10 = @(3, 3)
    11 = @(9, 9)
    12 = @(5, 5)
    13 = @(2, 2)
    14 = @(1, 1)
    15 = @(4, 4)
    16 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
